$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.658.58'
$ws.Range('D3').Value = '1.814.05'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''226.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('E6').Value = '  +3.65%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''38.46'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.95%  '
$ws.Range('E9').Value = '  -3.35%  '
$ws.Range('D10').Value = '''0.0682'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('D11').Value = '''0.0972'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('D12').Value = '2.076.29'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = '''11.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').Value = '1.823.39'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').Value = '''0.636'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').Value = '34.619.46'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = '''68.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').Value = '''244.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E20').Value = '  -2.57%  '
$ws.Range('D21').Value = '''11.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('E24').Value = '  +4.90%  '
$ws.Range('D25').Value = '''172.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = '''7.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').Value = '''17.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.11%  '
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '''3.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('D33').Value = '''1.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').Value = '''1.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').Value = '1.367.38'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('E36').Value = '  -3.65%  '
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''2.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.81%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0189'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').Value = '''1.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.40%  '
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''0.945'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''81.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '''14.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.62%  '
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '1.976.66'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('E48').Value = '  -3.77%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '''103.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').Value = '''49.27'
$ws.Range('D51').Style = 'Normal'
